$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows (Partida, LargoCm, AnchoCm, AltoCm, Prioridad, Remontable)
$data = @(
    ,@('(''SBCN25040629'', ''EBCN2517914'')', 80, 120, 102, 2, 0)
    ,@('(''SBCN25040359'', ''EBCN2517757'')', 80, 120, 71, 2, 0)
    ,@('(''SBCN25040609'', ''EBCN2517913'')', 120, 80, 68, 2, 0)
    ,@('(''SBCN25040331'', ''EBCN2517748'')', 78, 64, 108, 2, 0)
    ,@('(''SBCN25040313'', ''EBCN2517556'')', 89, 62, 90, 2, 0)
    ,@('(''SBCN25040330'', ''EBCN2517562'')', 61, 77, 102, 2, 0)
    ,@('(''SBCN25040108'', ''EBCN2517554'')', 61, 79, 88, 2, 0)
    ,@('(''SBCN25040628'', ''EBCN2517914'')', 120, 80, 39, 2, 0)
    ,@('(''SBCN25040604'', ''EBCN2517897'')', 79, 65, 64, 2, 0)
    ,@('(''SBCN25040360'', ''EBCN2517758'')', 80, 64, 51, 2, 0)
    ,@('(''SBCN25040317'', ''EBCN2517744'')', 36, 32, 25, 2, 1)
    ,@('(''SBCN25040358'', ''EBCN2517756'')', 50, 32, 22, 2, 1)
)

$firstRow = 2
$lastRow = $firstRow + $data.Length - 1

# Apply the same formatting used on the header cells (column A) to the new column A cells
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A" + $firstRow + ":A" + $lastRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $firstRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}
